$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 block: update semiEixo (C2) value
$ws.Range("C2").Value = 0.95

# Row 5 block: update lat (B5), semiEixo (C5), auxSemiEixo (D5) values
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 0.602
$ws.Range("D5").Value = 0.0559

# K10 holds a pasted-as-value copy of the (recalculated) K5 result
$ws.Range("K10").Value = 89.017641452899696

# Move the active selection from K10 to B2
$ws.Range("B2").Select()
